$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix row 5 "Name" cell: drop the trailing space in "Frankie Zayas " ---
$ws.Range("C5").Value = "Frankie Zayas"

# --- 2. Append a new form submission as row 10 ---
# First, clone row 8's look (row 8 is the previous "white stripe" row, the same
# banding row 10 will belong to) so fonts/fills/borders/number formats match.
$ws.Range("A8:AI8").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats

# Row 10 is now the last row of the table, so it needs the purple bottom rule
# (mirrors how row 9 looked before it stopped being the last row).
$ws.Range("C10:AI10").Borders.Item(9).Color = 6631236   # 442F65 bottom border

# --- 3. Fill in the new submission's values ---
$ws.Range("A10").Value = 44680.275181550925
$ws.Range("B10").Value = "zags66@yahoo.com"
$ws.Range("C10").Value = "Kraig Kuster"
$ws.Range("D10").Value = "Cam Ward -- Miami (Fla.) Sr QB"
$ws.Range("E10").Value = "Travis Hunter -- Colorado Jr ATH"
$ws.Range("F10").Value = "Abdul Carter -- Penn State Jr EDGE"
$ws.Range("G10").Value = "Will Campbell -- LSU Jr OT"
$ws.Range("H10").Value = "Ashton Jeanty -- Boise St. Jr RB"
$ws.Range("I10").Value = "Jalon Walker -- Georgia Jr EDGE"
$ws.Range("J10").Value = "Tyler Warren -- Penn State Sr TE"
$ws.Range("K10").Value = "Mykel Williams -- Georgia Jr EDGE"
$ws.Range("L10").Value = "Armand Membou -- Missouri Jr OT"
$ws.Range("M10").Value = "Kelvin Banks Jr. -- Texas Jr OT"
$ws.Range("N10").Value = "Walter Nolen -- Ole Miss Jr DL"
$ws.Range("O10").Value = "Matthew Golden -- Texas Jr WR"
$ws.Range("P10").Value = "Malaki Starks -- Georgia Jr S"
$ws.Range("Q10").Value = "Colston Loveland -- Michigan Jr TE"
$ws.Range("R10").Value = "Mike Green -- Marshall Soph EDGE"
$ws.Range("S10").Value = "Tetairoa McMillan -- Arizona Jr WR"
$ws.Range("T10").Value = "Shemar Stewart -- Texas A&M Jr EDGE"
$ws.Range("U10").Value = "Grey Zabel -- N. Dakota St. Sr IOL"
$ws.Range("V10").Value = "Jihaad Campbell -- Alabama Jr LB"
$ws.Range("W10").Value = "Omarion Hampton -- North Carolina Jr RB"
$ws.Range("X10").Value = "Shedeur Sanders -- Colorado Sr QB"
$ws.Range("Y10").Value = "Kenneth Grant -- Michigan Jr DL"
$ws.Range("Z10").Value = "Emeka Egbuka -- Ohio State Sr WR"
$ws.Range("AA10").Value = "Jaxson Dart -- Ole Miss Sr QB"
$ws.Range("AB10").Value = "Tyler Booker -- Alabama Jr IOL"
$ws.Range("AC10").Value = "Jahdae Barron -- Texas Sr CB"
$ws.Range("AD10").Value = "Shavon Revel Jr. -- East Carolina Sr CB"
$ws.Range("AE10").Value = "Derrick Harmon -- Oregon Jr DL"
$ws.Range("AF10").Value = "Donovan Ezeiruaku -- Boston College Sr EDGE"
$ws.Range("AG10").Value = "Nick Emmanwori -- South Carolina Jr S"
$ws.Range("AH10").Value = "Josh Simmons -- Ohio State Sr OT"
$ws.Range("AI10").Value = "Josh Conerly Jr. -- Oregon Jr OT"

# --- 4. Grow the table so the new row is part of Table_1 ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("C2:AI10"))
